$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 79
$ws1.Range("F9").Value = 8877
$ws1.Range("G11").Value = "已售罄"
$ws1.Range("F13").Value = 1005
$ws1.Range("F14").Value = 120
$ws1.Range("F19").Value = 71
$ws1.Range("F21").Value = 1079

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 79
$ws4.Range("F11").Value = 8877
$ws4.Range("G13").Value = "已售罄"
$ws4.Range("F15").Value = 1005
$ws4.Range("F16").Value = 120
$ws4.Range("F21").Value = 71
$ws4.Range("F23").Value = 1079
